$wb = $excel.ActiveWorkbook

# --- Rename the original sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Layer Centrality"

# --- Add the two cluster sheets, in order, after the last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "cluster_0"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "cluster_-1"

# --- Copy the header/label formatting (bold, centered, bordered) from the
#     original sheet onto the new sheets before filling in values, so the
#     pasted style id is reused instead of new styles being created. ---
$ws1.Range("B1").Copy()
$ws2.Range("B1:H1").PasteSpecial(-4122)
$ws3.Range("B1:H1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$ws2.Range("A2:A60").PasteSpecial(-4122)
$ws3.Range("A2:A5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws2.Cells.Item(1, 2).Value = "coauthor"
$ws2.Cells.Item(1, 3).Value = "facebook"
$ws2.Cells.Item(1, 4).Value = "leisure"
$ws2.Cells.Item(1, 5).Value = "lunch"
$ws2.Cells.Item(1, 6).Value = "work"
$ws2.Cells.Item(1, 7).Value = "shannon_entropy"
$ws2.Cells.Item(1, 8).Value = "cluster_class"
$ws2.Cells.Item(2, 1).Value = "U1"
$ws2.Cells.Item(2, 2).Value = 1.32
$ws2.Cells.Item(2, 3).Value = 21.21
$ws2.Cells.Item(2, 4).Value = 16.93
$ws2.Cells.Item(2, 5).Value = 21.47
$ws2.Cells.Item(2, 6).Value = 39.07
$ws2.Cells.Item(2, 7).Value = 2
$ws2.Cells.Item(2, 8).Value = 0
$ws2.Cells.Item(3, 1).Value = "U10"
$ws2.Cells.Item(3, 2).Value = 1.33
$ws2.Cells.Item(3, 3).Value = 23.81
$ws2.Cells.Item(3, 4).Value = 23.49
$ws2.Cells.Item(3, 5).Value = 18.11
$ws2.Cells.Item(3, 6).Value = 33.26
$ws2.Cells.Item(3, 7).Value = 2.04
$ws2.Cells.Item(3, 8).Value = 0
$ws2.Cells.Item(4, 1).Value = "U106"
$ws2.Cells.Item(4, 2).Value = 1.46
$ws2.Cells.Item(4, 3).Value = 30.35
$ws2.Cells.Item(4, 4).Value = 3.73
$ws2.Cells.Item(4, 5).Value = 30.96
$ws2.Cells.Item(4, 6).Value = 33.51
$ws2.Cells.Item(4, 7).Value = 1.84
$ws2.Cells.Item(4, 8).Value = 0
$ws2.Cells.Item(5, 1).Value = "U107"
$ws2.Cells.Item(5, 2).Value = 0.82
$ws2.Cells.Item(5, 3).Value = 25.49
$ws2.Cells.Item(5, 4).Value = 20.61
$ws2.Cells.Item(5, 5).Value = 19.75
$ws2.Cells.Item(5, 6).Value = 33.32
$ws2.Cells.Item(5, 7).Value = 2.02
$ws2.Cells.Item(5, 8).Value = 0
$ws2.Cells.Item(6, 1).Value = "U109"
$ws2.Cells.Item(6, 2).Value = 0.75
$ws2.Cells.Item(6, 3).Value = 17.72
$ws2.Cells.Item(6, 4).Value = 17.37
$ws2.Cells.Item(6, 5).Value = 37.93
$ws2.Cells.Item(6, 6).Value = 26.22
$ws2.Cells.Item(6, 7).Value = 1.97
$ws2.Cells.Item(6, 8).Value = 0
$ws2.Cells.Item(7, 1).Value = "U110"
$ws2.Cells.Item(7, 2).Value = 4.58
$ws2.Cells.Item(7, 3).Value = 19.91
$ws2.Cells.Item(7, 4).Value = 15.65
$ws2.Cells.Item(7, 5).Value = 23.91
$ws2.Cells.Item(7, 6).Value = 35.94
$ws2.Cells.Item(7, 7).Value = 2.11
$ws2.Cells.Item(7, 8).Value = 0
$ws2.Cells.Item(8, 1).Value = "U112"
$ws2.Cells.Item(8, 2).Value = 0.6899999999999999
$ws2.Cells.Item(8, 3).Value = 26.96
$ws2.Cells.Item(8, 4).Value = 2.76
$ws2.Cells.Item(8, 5).Value = 33.42
$ws2.Cells.Item(8, 6).Value = 36.17
$ws2.Cells.Item(8, 7).Value = 1.76
$ws2.Cells.Item(8, 8).Value = 0
$ws2.Cells.Item(9, 1).Value = "U113"
$ws2.Cells.Item(9, 2).Value = 0.93
$ws2.Cells.Item(9, 3).Value = 21.25
$ws2.Cells.Item(9, 4).Value = 18.14
$ws2.Cells.Item(9, 5).Value = 29.95
$ws2.Cells.Item(9, 6).Value = 29.72
$ws2.Cells.Item(9, 7).Value = 2.03
$ws2.Cells.Item(9, 8).Value = 0
$ws2.Cells.Item(10, 1).Value = "U118"
$ws2.Cells.Item(10, 2).Value = 1.51
$ws2.Cells.Item(10, 3).Value = 8.77
$ws2.Cells.Item(10, 4).Value = 2.79
$ws2.Cells.Item(10, 5).Value = 28.8
$ws2.Cells.Item(10, 6).Value = 58.13
$ws2.Cells.Item(10, 7).Value = 1.52
$ws2.Cells.Item(10, 8).Value = 0
$ws2.Cells.Item(11, 1).Value = "U123"
$ws2.Cells.Item(11, 2).Value = 0.64
$ws2.Cells.Item(11, 3).Value = 19.31
$ws2.Cells.Item(11, 4).Value = 2.45
$ws2.Cells.Item(11, 5).Value = 23.28
$ws2.Cells.Item(11, 6).Value = 54.32
$ws2.Cells.Item(11, 7).Value = 1.6
$ws2.Cells.Item(11, 8).Value = 0
$ws2.Cells.Item(12, 1).Value = "U124"
$ws2.Cells.Item(12, 2).Value = 0.88
$ws2.Cells.Item(12, 3).Value = 28.11
$ws2.Cells.Item(12, 4).Value = 15.75
$ws2.Cells.Item(12, 5).Value = 21.8
$ws2.Cells.Item(12, 6).Value = 33.46
$ws2.Cells.Item(12, 7).Value = 2
$ws2.Cells.Item(12, 8).Value = 0
$ws2.Cells.Item(13, 1).Value = "U126"
$ws2.Cells.Item(13, 2).Value = 0.36
$ws2.Cells.Item(13, 3).Value = 2.69
$ws2.Cells.Item(13, 4).Value = 24.9
$ws2.Cells.Item(13, 5).Value = 42.28
$ws2.Cells.Item(13, 6).Value = 29.76
$ws2.Cells.Item(13, 7).Value = 1.71
$ws2.Cells.Item(13, 8).Value = 0
$ws2.Cells.Item(14, 1).Value = "U13"
$ws2.Cells.Item(14, 2).Value = 1.54
$ws2.Cells.Item(14, 3).Value = 3.96
$ws2.Cells.Item(14, 4).Value = 20.07
$ws2.Cells.Item(14, 5).Value = 34.93
$ws2.Cells.Item(14, 6).Value = 39.5
$ws2.Cells.Item(14, 7).Value = 1.8
$ws2.Cells.Item(14, 8).Value = 0
$ws2.Cells.Item(15, 1).Value = "U130"
$ws2.Cells.Item(15, 2).Value = 4.24
$ws2.Cells.Item(15, 3).Value = 22.25
$ws2.Cells.Item(15, 4).Value = 3.82
$ws2.Cells.Item(15, 5).Value = 28.93
$ws2.Cells.Item(15, 6).Value = 40.77
$ws2.Cells.Item(15, 7).Value = 1.9
$ws2.Cells.Item(15, 8).Value = 0
$ws2.Cells.Item(16, 1).Value = "U134"
$ws2.Cells.Item(16, 2).Value = 2.54
$ws2.Cells.Item(16, 3).Value = 16.05
$ws2.Cells.Item(16, 4).Value = 3.98
$ws2.Cells.Item(16, 5).Value = 46.62
$ws2.Cells.Item(16, 6).Value = 30.8
$ws2.Cells.Item(16, 7).Value = 1.78
$ws2.Cells.Item(16, 8).Value = 0
$ws2.Cells.Item(17, 1).Value = "U138"
$ws2.Cells.Item(17, 2).Value = 3.26
$ws2.Cells.Item(17, 3).Value = 6.35
$ws2.Cells.Item(17, 4).Value = 21.92
$ws2.Cells.Item(17, 5).Value = 29.41
$ws2.Cells.Item(17, 6).Value = 39.07
$ws2.Cells.Item(17, 7).Value = 1.94
$ws2.Cells.Item(17, 8).Value = 0
$ws2.Cells.Item(18, 1).Value = "U14"
$ws2.Cells.Item(18, 2).Value = 0.52
$ws2.Cells.Item(18, 3).Value = 5.09
$ws2.Cells.Item(18, 4).Value = 20.16
$ws2.Cells.Item(18, 5).Value = 28.22
$ws2.Cells.Item(18, 6).Value = 46
$ws2.Cells.Item(18, 7).Value = 1.75
$ws2.Cells.Item(18, 8).Value = 0
$ws2.Cells.Item(19, 1).Value = "U141"
$ws2.Cells.Item(19, 2).Value = 0.35
$ws2.Cells.Item(19, 3).Value = 7.38
$ws2.Cells.Item(19, 4).Value = 10.91
$ws2.Cells.Item(19, 5).Value = 38.88
$ws2.Cells.Item(19, 6).Value = 42.47
$ws2.Cells.Item(19, 7).Value = 1.71
$ws2.Cells.Item(19, 8).Value = 0
$ws2.Cells.Item(20, 1).Value = "U142"
$ws2.Cells.Item(20, 2).Value = 1.46
$ws2.Cells.Item(20, 3).Value = 31.26
$ws2.Cells.Item(20, 4).Value = 17.17
$ws2.Cells.Item(20, 5).Value = 22.47
$ws2.Cells.Item(20, 6).Value = 27.64
$ws2.Cells.Item(20, 7).Value = 2.05
$ws2.Cells.Item(20, 8).Value = 0
$ws2.Cells.Item(21, 1).Value = "U17"
$ws2.Cells.Item(21, 2).Value = 0.49
$ws2.Cells.Item(21, 3).Value = 4.39
$ws2.Cells.Item(21, 4).Value = 26.32
$ws2.Cells.Item(21, 5).Value = 22.71
$ws2.Cells.Item(21, 6).Value = 46.09
$ws2.Cells.Item(21, 7).Value = 1.74
$ws2.Cells.Item(21, 8).Value = 0
$ws2.Cells.Item(22, 1).Value = "U18"
$ws2.Cells.Item(22, 2).Value = 3.41
$ws2.Cells.Item(22, 3).Value = 21.36
$ws2.Cells.Item(22, 4).Value = 15.92
$ws2.Cells.Item(22, 5).Value = 27.37
$ws2.Cells.Item(22, 6).Value = 31.94
$ws2.Cells.Item(22, 7).Value = 2.1
$ws2.Cells.Item(22, 8).Value = 0
$ws2.Cells.Item(23, 1).Value = "U19"
$ws2.Cells.Item(23, 2).Value = 0.39
$ws2.Cells.Item(23, 3).Value = 5.97
$ws2.Cells.Item(23, 4).Value = 17.12
$ws2.Cells.Item(23, 5).Value = 23.07
$ws2.Cells.Item(23, 6).Value = 53.45
$ws2.Cells.Item(23, 7).Value = 1.68
$ws2.Cells.Item(23, 8).Value = 0
$ws2.Cells.Item(24, 1).Value = "U21"
$ws2.Cells.Item(24, 2).Value = 0.79
$ws2.Cells.Item(24, 3).Value = 20.47
$ws2.Cells.Item(24, 4).Value = 14.94
$ws2.Cells.Item(24, 5).Value = 33.95
$ws2.Cells.Item(24, 6).Value = 29.85
$ws2.Cells.Item(24, 7).Value = 1.98
$ws2.Cells.Item(24, 8).Value = 0
$ws2.Cells.Item(25, 1).Value = "U22"
$ws2.Cells.Item(25, 2).Value = 1.51
$ws2.Cells.Item(25, 3).Value = 4.59
$ws2.Cells.Item(25, 4).Value = 15.35
$ws2.Cells.Item(25, 5).Value = 30.59
$ws2.Cells.Item(25, 6).Value = 47.96
$ws2.Cells.Item(25, 7).Value = 1.74
$ws2.Cells.Item(25, 8).Value = 0
$ws2.Cells.Item(26, 1).Value = "U23"
$ws2.Cells.Item(26, 2).Value = 0.52
$ws2.Cells.Item(26, 3).Value = 7.56
$ws2.Cells.Item(26, 4).Value = 24.74
$ws2.Cells.Item(26, 5).Value = 23.06
$ws2.Cells.Item(26, 6).Value = 44.12
$ws2.Cells.Item(26, 7).Value = 1.83
$ws2.Cells.Item(26, 8).Value = 0
$ws2.Cells.Item(27, 1).Value = "U26"
$ws2.Cells.Item(27, 2).Value = 1.14
$ws2.Cells.Item(27, 3).Value = 3.16
$ws2.Cells.Item(27, 4).Value = 12.64
$ws2.Cells.Item(27, 5).Value = 23.65
$ws2.Cells.Item(27, 6).Value = 59.4
$ws2.Cells.Item(27, 7).Value = 1.55
$ws2.Cells.Item(27, 8).Value = 0
$ws2.Cells.Item(28, 1).Value = "U29"
$ws2.Cells.Item(28, 2).Value = 1.45
$ws2.Cells.Item(28, 3).Value = 24.94
$ws2.Cells.Item(28, 4).Value = 19.05
$ws2.Cells.Item(28, 5).Value = 24.96
$ws2.Cells.Item(28, 6).Value = 29.6
$ws2.Cells.Item(28, 7).Value = 2.06
$ws2.Cells.Item(28, 8).Value = 0
$ws2.Cells.Item(29, 1).Value = "U3"
$ws2.Cells.Item(29, 2).Value = 0.72
$ws2.Cells.Item(29, 3).Value = 24.23
$ws2.Cells.Item(29, 4).Value = 14.91
$ws2.Cells.Item(29, 5).Value = 28.06
$ws2.Cells.Item(29, 6).Value = 32.09
$ws2.Cells.Item(29, 7).Value = 2
$ws2.Cells.Item(29, 8).Value = 0
$ws2.Cells.Item(30, 1).Value = "U32"
$ws2.Cells.Item(30, 2).Value = 1.44
$ws2.Cells.Item(30, 3).Value = 28.16
$ws2.Cells.Item(30, 4).Value = 10.95
$ws2.Cells.Item(30, 5).Value = 23.33
$ws2.Cells.Item(30, 6).Value = 36.11
$ws2.Cells.Item(30, 7).Value = 1.97
$ws2.Cells.Item(30, 8).Value = 0
$ws2.Cells.Item(31, 1).Value = "U33"
$ws2.Cells.Item(31, 2).Value = 0
$ws2.Cells.Item(31, 3).Value = 3.91
$ws2.Cells.Item(31, 4).Value = 0.43
$ws2.Cells.Item(31, 5).Value = 34.91
$ws2.Cells.Item(31, 6).Value = 60.76
$ws2.Cells.Item(31, 7).Value = 1.18
$ws2.Cells.Item(31, 8).Value = 0
$ws2.Cells.Item(32, 1).Value = "U37"
$ws2.Cells.Item(32, 2).Value = 0.33
$ws2.Cells.Item(32, 3).Value = 4.04
$ws2.Cells.Item(32, 4).Value = 23.35
$ws2.Cells.Item(32, 5).Value = 40.91
$ws2.Cells.Item(32, 6).Value = 31.36
$ws2.Cells.Item(32, 7).Value = 1.76
$ws2.Cells.Item(32, 8).Value = 0
$ws2.Cells.Item(33, 1).Value = "U4"
$ws2.Cells.Item(33, 2).Value = 0.77
$ws2.Cells.Item(33, 3).Value = 22.65
$ws2.Cells.Item(33, 4).Value = 7.54
$ws2.Cells.Item(33, 5).Value = 29.17
$ws2.Cells.Item(33, 6).Value = 39.87
$ws2.Cells.Item(33, 7).Value = 1.87
$ws2.Cells.Item(33, 8).Value = 0
$ws2.Cells.Item(34, 1).Value = "U41"
$ws2.Cells.Item(34, 2).Value = 1.55
$ws2.Cells.Item(34, 3).Value = 12.66
$ws2.Cells.Item(34, 4).Value = 4.4
$ws2.Cells.Item(34, 5).Value = 39.51
$ws2.Cells.Item(34, 6).Value = 41.88
$ws2.Cells.Item(34, 7).Value = 1.72
$ws2.Cells.Item(34, 8).Value = 0
$ws2.Cells.Item(35, 1).Value = "U42"
$ws2.Cells.Item(35, 2).Value = 1
$ws2.Cells.Item(35, 3).Value = 18.65
$ws2.Cells.Item(35, 4).Value = 18.28
$ws2.Cells.Item(35, 5).Value = 30.79
$ws2.Cells.Item(35, 6).Value = 31.28
$ws2.Cells.Item(35, 7).Value = 2.01
$ws2.Cells.Item(35, 8).Value = 0
$ws2.Cells.Item(36, 1).Value = "U47"
$ws2.Cells.Item(36, 2).Value = 3.41
$ws2.Cells.Item(36, 3).Value = 28.59
$ws2.Cells.Item(36, 4).Value = 12.37
$ws2.Cells.Item(36, 5).Value = 26.57
$ws2.Cells.Item(36, 6).Value = 29.06
$ws2.Cells.Item(36, 7).Value = 2.08
$ws2.Cells.Item(36, 8).Value = 0
$ws2.Cells.Item(37, 1).Value = "U48"
$ws2.Cells.Item(37, 2).Value = 1.37
$ws2.Cells.Item(37, 3).Value = 4.46
$ws2.Cells.Item(37, 4).Value = 5.16
$ws2.Cells.Item(37, 5).Value = 45.98
$ws2.Cells.Item(37, 6).Value = 43.02
$ws2.Cells.Item(37, 7).Value = 1.54
$ws2.Cells.Item(37, 8).Value = 0
$ws2.Cells.Item(38, 1).Value = "U49"
$ws2.Cells.Item(38, 2).Value = 0.87
$ws2.Cells.Item(38, 3).Value = 17.91
$ws2.Cells.Item(38, 4).Value = 12.9
$ws2.Cells.Item(38, 5).Value = 28.19
$ws2.Cells.Item(38, 6).Value = 40.13
$ws2.Cells.Item(38, 7).Value = 1.93
$ws2.Cells.Item(38, 8).Value = 0
$ws2.Cells.Item(39, 1).Value = "U53"
$ws2.Cells.Item(39, 2).Value = 5.3
$ws2.Cells.Item(39, 3).Value = 10.28
$ws2.Cells.Item(39, 4).Value = 20.49
$ws2.Cells.Item(39, 5).Value = 32.03
$ws2.Cells.Item(39, 6).Value = 31.9
$ws2.Cells.Item(39, 7).Value = 2.08
$ws2.Cells.Item(39, 8).Value = 0
$ws2.Cells.Item(40, 1).Value = "U54"
$ws2.Cells.Item(40, 2).Value = 0.7
$ws2.Cells.Item(40, 3).Value = 18.75
$ws2.Cells.Item(40, 4).Value = 20.6
$ws2.Cells.Item(40, 5).Value = 26.52
$ws2.Cells.Item(40, 6).Value = 33.42
$ws2.Cells.Item(40, 7).Value = 2.01
$ws2.Cells.Item(40, 8).Value = 0
$ws2.Cells.Item(41, 1).Value = "U59"
$ws2.Cells.Item(41, 2).Value = 0.96
$ws2.Cells.Item(41, 3).Value = 21.83
$ws2.Cells.Item(41, 4).Value = 17.73
$ws2.Cells.Item(41, 5).Value = 31.96
$ws2.Cells.Item(41, 6).Value = 27.52
$ws2.Cells.Item(41, 7).Value = 2.02
$ws2.Cells.Item(41, 8).Value = 0
$ws2.Cells.Item(42, 1).Value = "U6"
$ws2.Cells.Item(42, 2).Value = 0.72
$ws2.Cells.Item(42, 3).Value = 19.75
$ws2.Cells.Item(42, 4).Value = 14.3
$ws2.Cells.Item(42, 5).Value = 33
$ws2.Cells.Item(42, 6).Value = 32.24
$ws2.Cells.Item(42, 7).Value = 1.97
$ws2.Cells.Item(42, 8).Value = 0
$ws2.Cells.Item(43, 1).Value = "U62"
$ws2.Cells.Item(43, 2).Value = 0.33
$ws2.Cells.Item(43, 3).Value = 3.31
$ws2.Cells.Item(43, 4).Value = 13.89
$ws2.Cells.Item(43, 5).Value = 39.26
$ws2.Cells.Item(43, 6).Value = 43.21
$ws2.Cells.Item(43, 7).Value = 1.64
$ws2.Cells.Item(43, 8).Value = 0
$ws2.Cells.Item(44, 1).Value = "U63"
$ws2.Cells.Item(44, 2).Value = 0
$ws2.Cells.Item(44, 3).Value = 3.95
$ws2.Cells.Item(44, 4).Value = 1.5
$ws2.Cells.Item(44, 5).Value = 42.25
$ws2.Cells.Item(44, 6).Value = 52.3
$ws2.Cells.Item(44, 7).Value = 1.29
$ws2.Cells.Item(44, 8).Value = 0
$ws2.Cells.Item(45, 1).Value = "U65"
$ws2.Cells.Item(45, 2).Value = 0.86
$ws2.Cells.Item(45, 3).Value = 21.87
$ws2.Cells.Item(45, 4).Value = 23.79
$ws2.Cells.Item(45, 5).Value = 26.24
$ws2.Cells.Item(45, 6).Value = 27.24
$ws2.Cells.Item(45, 7).Value = 2.05
$ws2.Cells.Item(45, 8).Value = 0
$ws2.Cells.Item(46, 1).Value = "U67"
$ws2.Cells.Item(46, 2).Value = 0.79
$ws2.Cells.Item(46, 3).Value = 21.22
$ws2.Cells.Item(46, 4).Value = 10.02
$ws2.Cells.Item(46, 5).Value = 25.82
$ws2.Cells.Item(46, 6).Value = 42.15
$ws2.Cells.Item(46, 7).Value = 1.89
$ws2.Cells.Item(46, 8).Value = 0
$ws2.Cells.Item(47, 1).Value = "U68"
$ws2.Cells.Item(47, 2).Value = 1.42
$ws2.Cells.Item(47, 3).Value = 9.74
$ws2.Cells.Item(47, 4).Value = 15.58
$ws2.Cells.Item(47, 5).Value = 35.5
$ws2.Cells.Item(47, 6).Value = 37.76
$ws2.Cells.Item(47, 7).Value = 1.89
$ws2.Cells.Item(47, 8).Value = 0
$ws2.Cells.Item(48, 1).Value = "U69"
$ws2.Cells.Item(48, 2).Value = 0.74
$ws2.Cells.Item(48, 3).Value = 23.53
$ws2.Cells.Item(48, 4).Value = 15.66
$ws2.Cells.Item(48, 5).Value = 30.33
$ws2.Cells.Item(48, 6).Value = 29.75
$ws2.Cells.Item(48, 7).Value = 2
$ws2.Cells.Item(48, 8).Value = 0
$ws2.Cells.Item(49, 1).Value = "U71"
$ws2.Cells.Item(49, 2).Value = 0.74
$ws2.Cells.Item(49, 3).Value = 20.79
$ws2.Cells.Item(49, 4).Value = 2.5
$ws2.Cells.Item(49, 5).Value = 26.62
$ws2.Cells.Item(49, 6).Value = 49.35
$ws2.Cells.Item(49, 7).Value = 1.67
$ws2.Cells.Item(49, 8).Value = 0
$ws2.Cells.Item(50, 1).Value = "U72"
$ws2.Cells.Item(50, 2).Value = 3.13
$ws2.Cells.Item(50, 3).Value = 6.84
$ws2.Cells.Item(50, 4).Value = 19.57
$ws2.Cells.Item(50, 5).Value = 37.87
$ws2.Cells.Item(50, 6).Value = 32.58
$ws2.Cells.Item(50, 7).Value = 1.94
$ws2.Cells.Item(50, 8).Value = 0
$ws2.Cells.Item(51, 1).Value = "U73"
$ws2.Cells.Item(51, 2).Value = 0.5600000000000001
$ws2.Cells.Item(51, 3).Value = 6.75
$ws2.Cells.Item(51, 4).Value = 32.08
$ws2.Cells.Item(51, 5).Value = 25.57
$ws2.Cells.Item(51, 6).Value = 35.04
$ws2.Cells.Item(51, 7).Value = 1.86
$ws2.Cells.Item(51, 8).Value = 0
$ws2.Cells.Item(52, 1).Value = "U76"
$ws2.Cells.Item(52, 2).Value = 3.1
$ws2.Cells.Item(52, 3).Value = 20.82
$ws2.Cells.Item(52, 4).Value = 17.84
$ws2.Cells.Item(52, 5).Value = 27.14
$ws2.Cells.Item(52, 6).Value = 31.11
$ws2.Cells.Item(52, 7).Value = 2.11
$ws2.Cells.Item(52, 8).Value = 0
$ws2.Cells.Item(53, 1).Value = "U79"
$ws2.Cells.Item(53, 2).Value = 0.79
$ws2.Cells.Item(53, 3).Value = 23.02
$ws2.Cells.Item(53, 4).Value = 17.64
$ws2.Cells.Item(53, 5).Value = 28.77
$ws2.Cells.Item(53, 6).Value = 29.78
$ws2.Cells.Item(53, 7).Value = 2.02
$ws2.Cells.Item(53, 8).Value = 0
$ws2.Cells.Item(54, 1).Value = "U86"
$ws2.Cells.Item(54, 2).Value = 0.74
$ws2.Cells.Item(54, 3).Value = 6.53
$ws2.Cells.Item(54, 4).Value = 0.64
$ws2.Cells.Item(54, 5).Value = 40.7
$ws2.Cells.Item(54, 6).Value = 51.38
$ws2.Cells.Item(54, 7).Value = 1.38
$ws2.Cells.Item(54, 8).Value = 0
$ws2.Cells.Item(55, 1).Value = "U90"
$ws2.Cells.Item(55, 2).Value = 0.36
$ws2.Cells.Item(55, 3).Value = 3.05
$ws2.Cells.Item(55, 4).Value = 25.17
$ws2.Cells.Item(55, 5).Value = 33.85
$ws2.Cells.Item(55, 6).Value = 37.58
$ws2.Cells.Item(55, 7).Value = 1.74
$ws2.Cells.Item(55, 8).Value = 0
$ws2.Cells.Item(56, 1).Value = "U91"
$ws2.Cells.Item(56, 2).Value = 2.75
$ws2.Cells.Item(56, 3).Value = 27.55
$ws2.Cells.Item(56, 4).Value = 24.51
$ws2.Cells.Item(56, 5).Value = 20.97
$ws2.Cells.Item(56, 6).Value = 24.21
$ws2.Cells.Item(56, 7).Value = 2.12
$ws2.Cells.Item(56, 8).Value = 0
$ws2.Cells.Item(57, 1).Value = "U92"
$ws2.Cells.Item(57, 2).Value = 0.02
$ws2.Cells.Item(57, 3).Value = 6.32
$ws2.Cells.Item(57, 4).Value = 0.72
$ws2.Cells.Item(57, 5).Value = 49.07
$ws2.Cells.Item(57, 6).Value = 43.87
$ws2.Cells.Item(57, 7).Value = 1.33
$ws2.Cells.Item(57, 8).Value = 0
$ws2.Cells.Item(58, 1).Value = "U97"
$ws2.Cells.Item(58, 2).Value = 3.84
$ws2.Cells.Item(58, 3).Value = 3.9
$ws2.Cells.Item(58, 4).Value = 3.16
$ws2.Cells.Item(58, 5).Value = 39.92
$ws2.Cells.Item(58, 6).Value = 49.19
$ws2.Cells.Item(58, 7).Value = 1.55
$ws2.Cells.Item(58, 8).Value = 0
$ws2.Cells.Item(59, 1).Value = "U99"
$ws2.Cells.Item(59, 2).Value = 3.97
$ws2.Cells.Item(59, 3).Value = 5.78
$ws2.Cells.Item(59, 4).Value = 17.92
$ws2.Cells.Item(59, 5).Value = 31.64
$ws2.Cells.Item(59, 6).Value = 40.69
$ws2.Cells.Item(59, 7).Value = 1.92
$ws2.Cells.Item(59, 8).Value = 0
$ws2.Cells.Item(60, 1).Value = "mean"
$ws2.Cells.Item(60, 2).Value = 1.416551724137931
$ws2.Cells.Item(60, 3).Value = 15.19310344827587
$ws2.Cells.Item(60, 4).Value = 14.38413793103448
$ws2.Cells.Item(60, 5).Value = 30.7401724137931
$ws2.Cells.Item(60, 6).Value = 38.265
$ws2.Cells.Item(60, 7).Value = 1.84051724137931
$ws2.Cells.Item(60, 8).Value = 0

$ws3.Cells.Item(1, 2).Value = "coauthor"
$ws3.Cells.Item(1, 3).Value = "facebook"
$ws3.Cells.Item(1, 4).Value = "leisure"
$ws3.Cells.Item(1, 5).Value = "lunch"
$ws3.Cells.Item(1, 6).Value = "work"
$ws3.Cells.Item(1, 7).Value = "shannon_entropy"
$ws3.Cells.Item(1, 8).Value = "cluster_class"
$ws3.Cells.Item(2, 1).Value = "U102"
$ws3.Cells.Item(2, 2).Value = 0
$ws3.Cells.Item(2, 3).Value = 4.71
$ws3.Cells.Item(2, 4).Value = 0.65
$ws3.Cells.Item(2, 5).Value = 58.79
$ws3.Cells.Item(2, 6).Value = 35.84
$ws3.Cells.Item(2, 7).Value = 1.24
$ws3.Cells.Item(2, 8).Value = -1
$ws3.Cells.Item(3, 1).Value = "U139"
$ws3.Cells.Item(3, 2).Value = 0
$ws3.Cells.Item(3, 3).Value = 0.9
$ws3.Cells.Item(3, 4).Value = 0.54
$ws3.Cells.Item(3, 5).Value = 16.58
$ws3.Cells.Item(3, 6).Value = 81.98
$ws3.Cells.Item(3, 7).Value = 0.77
$ws3.Cells.Item(3, 8).Value = -1
$ws3.Cells.Item(4, 1).Value = "U140"
$ws3.Cells.Item(4, 2).Value = 0
$ws3.Cells.Item(4, 3).Value = 0
$ws3.Cells.Item(4, 4).Value = 0
$ws3.Cells.Item(4, 5).Value = 4.85
$ws3.Cells.Item(4, 6).Value = 95.15000000000001
$ws3.Cells.Item(4, 7).Value = 0.28
$ws3.Cells.Item(4, 8).Value = -1
$ws3.Cells.Item(5, 1).Value = "mean"
$ws3.Cells.Item(5, 2).Value = 0
$ws3.Cells.Item(5, 3).Value = 1.87
$ws3.Cells.Item(5, 4).Value = 0.3966666666666667
$ws3.Cells.Item(5, 5).Value = 26.74
$ws3.Cells.Item(5, 6).Value = 70.99000000000001
$ws3.Cells.Item(5, 7).Value = 0.7633333333333333
$ws3.Cells.Item(5, 8).Value = -1
# --- Conditional formatting: four graduated "heat" thresholds, same
#     thresholds/colors as the original LayerCentrality sheet. ---
$cf2 = $ws2.Range("B2:F59")
$r = $cf2.FormatConditions.Add(1, 7, "75")
$r.Interior.Color = 2570461
$r = $cf2.FormatConditions.Add(1, 7, "50")
$r.Interior.Color = 1157359
$r = $cf2.FormatConditions.Add(1, 7, "30")
$r.Interior.Color = 324599
$r = $cf2.FormatConditions.Add(1, 7, "0")
$r.Interior.Color = 3065900

$cf3 = $ws3.Range("B2:F4")
$r = $cf3.FormatConditions.Add(1, 7, "75")
$r.Interior.Color = 2570461
$r = $cf3.FormatConditions.Add(1, 7, "50")
$r.Interior.Color = 1157359
$r = $cf3.FormatConditions.Add(1, 7, "30")
$r.Interior.Color = 324599
$r = $cf3.FormatConditions.Add(1, 7, "0")
$r.Interior.Color = 3065900

# --- Restore the original sheet as the active tab ---
$ws1.Activate()
